$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2 through 447) holds the "Förändrad" (changed) date as a
# serial date number. Every row currently stores 46082 (2026-03-01) and
# needs to be bumped by one day to 46083 (2026-03-02).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 447) { $lastRow = 447 }

$ws.Range("C2:C$lastRow").Value = 46083
